$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on the "展览" sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 132
$wsExhibit.Range("F3").Value = 34

# Update the same values on the "全部类型" sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 132
$wsAll.Range("F3").Value = 34
